$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB) {
    $rangeA = $ws.Range("B" + $rowA + ":AD" + $rowA)
    $rangeB = $ws.Range("B" + $rowB + ":AD" + $rowB)
    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()
    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

Swap-Rows $ws 117 118
Swap-Rows $ws 151 152
Swap-Rows $ws 174 175
Swap-Rows $ws 180 181
Swap-Rows $ws 184 185
